$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 75: the "Scalpel Accuracy" label + value moved from C75/D75 to E75/F75
$ws.Range("C75").Value = $null
$ws.Range("D75").Value = $null
$ws.Range("E75").Value = "Scalpel Accuracy:"
$ws.Range("F75").Value = 91.67

# Row 76: fix label text
$ws.Range("E76").Value = "Accuracy vs PyType"
